$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "re_profiles" sheet: rename the commodity "elc_win-ITA" -> "elc_won-ITA"
#    (shared string used across the whole K4:K123 commodity column)
# ---------------------------------------------------------------------------
$wsProfiles = $wb.Worksheets.Item("re_profiles")
$wsProfiles.Range("K4:K123").Value = "elc_won-ITA"

# ---------------------------------------------------------------------------
# 2) "re_profiles" sheet: the small hydro timeslice/value table (rows 4-9)
#    got its (season, value) pairs shuffled into a new row order.
# ---------------------------------------------------------------------------
$wsProfiles.Range("M4").Value = "S5"
$wsProfiles.Range("N4").Value = 0.093706617306802395

$wsProfiles.Range("M5").Value = "S6"
$wsProfiles.Range("N5").Value = 0.25116844053678855

$wsProfiles.Range("M6").Value = "S3"
$wsProfiles.Range("N6").Value = 0.37937297547431742

$wsProfiles.Range("M7").Value = "S1"
$wsProfiles.Range("N7").Value = 0.14770939379916706

$wsProfiles.Range("M8").Value = "S2"
$wsProfiles.Range("N8").Value = 0.075520592318371119

$wsProfiles.Range("M9").Value = "S4"
$wsProfiles.Range("N9").Value = 0.25252198056455344

# ---------------------------------------------------------------------------
# 3) "ev_charging_uc" sheet: the two comma-separated timeslice lists (C13 and
#    C14, surfaced via formulas in G7 =C14 and G8 =C13) had their comma
#    separated members reshuffled (same membership, new order).
# ---------------------------------------------------------------------------
$wsEv = $wb.Worksheets.Item("ev_charging_uc")
$wsEv.Range("C13").Value = "S1aH3,S6c1209h11,S6d1216h11,S6d1216h12,S6d1216h18,S4aH2,S6c1209h09,S6c1209h10,S1b0202h09,S6aH5,S1aH2,S1aH7,S1b0202h17,S2aH3,S3aH6,S4aH7,S5aH6,S5aH7,S6aH6,S6c1209h13,S6c1209h16,S1b0202h13,S1b0202h16,S5aH4,S6aH3,S6aH4,S6c1209h14,S6c1209h17,S1b0202h10,S1b0202h18,S2aH7,S3aH4,S4aH3,S4aH4,S6c1209h15,S6d1216h09,S6d1216h13,S6d1216h16,S1aH5,S1b0202h07,S2aH4,S2aH5,S5aH2,S5aH5,S6aH2,S1aH4,S1b0202h12,S2aH2,S4aH6,S6c1209h12,S6d1216h07,S6d1216h08,S6d1216h17,S1b0202h08,S1b0202h11,S3aH2,S6d1216h15,S1aH6,S6d1216h14,S3aH7,S6c1209h07,S1b0202h14,S2aH6,S3aH5,S5aH3,S1b0202h15,S3aH3,S4aH5,S6aH7,S6c1209h08,S6c1209h18,S6d1216h10"
$wsEv.Range("C14").Value = "S1aH1,S1b0202h20,S1b0202h23,S6d1216h24,S1b0202h03,S1b0202h22,S2aH1,S5aH8,S6c1209h02,S1b0202h04,S4aH8,S6c1209h06,S6c1209h24,S6c1209h20,S6d1216h05,S6d1216h01,S6d1216h20,S1b0202h05,S1b0202h06,S6c1209h23,S6d1216h02,S6d1216h23,S4aH1,S6aH1,S6d1216h04,S6d1216h19,S1b0202h24,S3aH8,S6c1209h04,S6aH8,S1b0202h02,S1b0202h21,S2aH8,S6c1209h03,S6c1209h05,S1aH8,S5aH1,S6c1209h19,S6c1209h22,S6d1216h03,S6d1216h06,S1b0202h01,S6d1216h21,S6c1209h21,S6d1216h22,S1b0202h19,S3aH1,S6c1209h01"

# Recalculate so the G7/G8 cached formula results (=C14 / =C13) pick up the
# new text.
$excel.Calculate()
